$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. The "Förändrad" (last-changed) column C was bumped from 45182 to 45184
#    for every existing data row (rows 2-538).
$ws.Range("C2:C538").Value = 45184

# 2. Row 538 picks up an explicit row height (it was the one row missing it
#    while every other row already had ht="15" customHeight="1").
$ws.Rows.Item(538).RowHeight = 15

# 3. Append a brand new record row (539) for notice "A 42874-2023".
$ws.Range("B539").NumberFormat = "YYYY-MM-DD"
$ws.Range("C539").NumberFormat = "YYYY-MM-DD"
$ws.Range("R539").WrapText = $true

$ws.Range("A539").Value = "A 42874-2023"
$ws.Range("B539").Value = 45182
$ws.Range("C539").Value = 45184
$ws.Range("D539").Value = "DALARNAS LÄN"
$ws.Range("E539").Value = "ÄLVDALEN"
$ws.Range("F539").Value = "Sveaskog"
$ws.Range("G539").Value = 2.7
$ws.Range("H539").Value = 0
$ws.Range("I539").Value = 0
$ws.Range("J539").Value = 0
$ws.Range("K539").Value = 0
$ws.Range("L539").Value = 0
$ws.Range("M539").Value = 0
$ws.Range("N539").Value = 0
$ws.Range("O539").Value = 0
$ws.Range("P539").Value = 0
$ws.Range("Q539").Value = 0
$ws.Range("R539").Value = ""
